$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data rows 2-6 (header row 1 unchanged)
# Columns: A trialTrain, B x_fixStart, C y_fixStart, D x_corrSteps, E y_corrSteps,
#          F x_nrSteps, G y_nrSteps, H alienID, I praclen, J version

$data = @(
    @(1, 7, 9, 2, 8, -5, -1, 12, 5),
    @(2, 5, 9, 1, 7, -4, -2, 23, 5),
    @(3, 9, 5, 8, 0, -1, -5, 56, 5),
    @(4, 8, 8, 5, 5, -3, -3, 34, 5),
    @(5, 6, 6, 4, 2, -2, -4, 45, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}
